$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Id values between row 2 and row 3
$ws.Range("A2").Value = 111659700
$ws.Range("A3").Value = 111659699

# Swap the coordinate values (Ost/Q and Nord/R) between row 2 and row 3
$ws.Range("Q2").Value = 799972.0195060072
$ws.Range("R2").Value = 7239766.02062137

$ws.Range("Q3").Value = 800047.4485974194
$ws.Range("R3").Value = 7239832.989497012
